# Table 1 update:
# Remove "Visits per year (median, IQR)" (baseline) row and replace it with a
# new "Visits per year during follow-up (median, IQR)" row placed after the
# two "Follow-up in days ..." rows (i.e. appended at the end of the table).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the old "Visits per year (median, IQR)" row (row 23: baseline visits/year).
$ws.Rows.Item(23).Delete()

# Append the new "Visits per year during follow-up (median, IQR)" row at the
# end of the table (now row 25, after the two follow-up-days rows).
$ws.Range("A25").Value = "Visits per year during follow-up (median, IQR)"
$ws.Range("B25").Value = "3 (2, 5)"
$ws.Range("C25").Value = "4 (3, 7)"
